# Update "想去人数" (column F) figures across all four sheets to match the
# freshly generated gh-pages data snapshot (commit 456a3b4).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(3, 6).Value = 10129
$ws.Cells.Item(5, 6).Value = 89
$ws.Cells.Item(7, 6).Value = 0
$ws.Cells.Item(8, 6).Value = 647
$ws.Cells.Item(9, 6).Value = 0
$ws.Cells.Item(10, 6).Value = 0
$ws.Cells.Item(11, 6).Value = 11506
$ws.Cells.Item(12, 6).Value = 12292
$ws.Cells.Item(13, 6).Value = 1285
$ws.Cells.Item(14, 6).Value = 0
$ws.Cells.Item(15, 6).Value = 5286
$ws.Cells.Item(16, 6).Value = 0
$ws.Cells.Item(21, 6).Value = 0
$ws.Cells.Item(22, 6).Value = 316
$ws.Cells.Item(23, 6).Value = 1967
$ws.Cells.Item(24, 6).Value = 978
$ws.Cells.Item(25, 6).Value = 1437
$ws.Cells.Item(26, 6).Value = 0
$ws.Cells.Item(27, 6).Value = 9
$ws.Cells.Item(28, 6).Value = 2132
$ws.Cells.Item(29, 6).Value = 477
$ws.Cells.Item(30, 6).Value = 695
$ws.Cells.Item(31, 6).Value = 0
$ws.Cells.Item(32, 6).Value = 231
$ws.Cells.Item(33, 6).Value = 1953
$ws.Cells.Item(34, 6).Value = 103
$ws.Cells.Item(35, 6).Value = 1438
$ws.Cells.Item(36, 6).Value = 130
$ws.Cells.Item(37, 6).Value = 966
$ws.Cells.Item(38, 6).Value = 92
$ws.Cells.Item(39, 6).Value = 84
$ws.Cells.Item(40, 6).Value = 0
$ws.Cells.Item(42, 6).Value = 111
$ws.Cells.Item(43, 6).Value = 572
$ws.Cells.Item(45, 6).Value = 0
$ws.Cells.Item(46, 6).Value = 909
$ws.Cells.Item(47, 6).Value = 270
$ws.Cells.Item(49, 6).Value = 4270
$ws.Cells.Item(50, 6).Value = 0

$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(4, 6).Value = 0
$ws.Cells.Item(5, 6).Value = 0
$ws.Cells.Item(6, 6).Value = 12
$ws.Cells.Item(8, 6).Value = 33
$ws.Cells.Item(9, 6).Value = 0
$ws.Cells.Item(15, 6).Value = 5
$ws.Cells.Item(16, 6).Value = 5
$ws.Cells.Item(23, 6).Value = 73
$ws.Cells.Item(25, 6).Value = 108
$ws.Cells.Item(26, 6).Value = 47
$ws.Cells.Item(29, 6).Value = 4

$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(2, 6).Value = 0

$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(2, 6).Value = 502
$ws.Cells.Item(3, 6).Value = 10129
$ws.Cells.Item(5, 6).Value = 0
$ws.Cells.Item(6, 6).Value = 14
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(10, 6).Value = 11506
$ws.Cells.Item(11, 6).Value = 0
$ws.Cells.Item(13, 6).Value = 1285
$ws.Cells.Item(14, 6).Value = 1254
$ws.Cells.Item(17, 6).Value = 0
$ws.Cells.Item(18, 6).Value = 58
$ws.Cells.Item(19, 6).Value = 0
$ws.Cells.Item(20, 6).Value = 0
$ws.Cells.Item(23, 6).Value = 316
$ws.Cells.Item(24, 6).Value = 1967
$ws.Cells.Item(25, 6).Value = 978
$ws.Cells.Item(27, 6).Value = 477
$ws.Cells.Item(28, 6).Value = 695
$ws.Cells.Item(29, 6).Value = 2891
$ws.Cells.Item(31, 6).Value = 0
$ws.Cells.Item(32, 6).Value = 103
$ws.Cells.Item(34, 6).Value = 1438
$ws.Cells.Item(35, 6).Value = 0
$ws.Cells.Item(36, 6).Value = 12
$ws.Cells.Item(37, 6).Value = 5
$ws.Cells.Item(39, 6).Value = 0
$ws.Cells.Item(40, 6).Value = 92
$ws.Cells.Item(41, 6).Value = 0
$ws.Cells.Item(42, 6).Value = 0
$ws.Cells.Item(43, 6).Value = 0
$ws.Cells.Item(44, 6).Value = 0
$ws.Cells.Item(45, 6).Value = 0
$ws.Cells.Item(46, 6).Value = 0
$ws.Cells.Item(47, 6).Value = 0
$ws.Cells.Item(49, 6).Value = 140
$ws.Cells.Item(50, 6).Value = 0
